# Updated symbol list on Fri Jan  6 19:21:23 UTC 2023 with GitHub Actions
#
# Refresh the "Price" (column D) and "Volume(1h)" (column E) figures for the
# crypto-exchange-token rows in the sheet to the latest scraped values.
# Values are written with a leading apostrophe so Excel keeps them as literal
# text (matching how the sheet already stores these figures) instead of
# re-interpreting "258.57" as a number or "0.46%" as a percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'258.57"
$ws.Range("E2").Value = "'0.46%"
$ws.Range("D3").Value = "'26.88"
$ws.Range("E3").Value = "'-2.02%"
$ws.Range("D4").Value = "'4.676"
$ws.Range("E4").Value = "'2.06%"
$ws.Range("D5").Value = "'0.05998"
$ws.Range("E5").Value = "'1.87%"
$ws.Range("D6").Value = "'6.666"
$ws.Range("E6").Value = "'0.56%"
$ws.Range("D7").Value = "'0.8577"
$ws.Range("E7").Value = "'-0.04%"
$ws.Range("D8").Value = "'0.9251"
$ws.Range("E8").Value = "'0.11%"
$ws.Range("D9").Value = "'0.1391"
$ws.Range("E9").Value = "'-1.47%"
$ws.Range("D10").Value = "'0.04546"
$ws.Range("E10").Value = "'27.43%"
$ws.Range("D11").Value = "'0.07034"
$ws.Range("E11").Value = "'-0.71%"
$ws.Range("D12").Value = "'0.03118"
$ws.Range("E12").Value = "'-3.37%"
$ws.Range("D13").Value = "'0.09130"
$ws.Range("E13").Value = "'-0.58%"
$ws.Range("D14").Value = "'0.001524"
$ws.Range("E14").Value = "'-0.96%"
$ws.Range("E15").Value = "'-0.32%"
$ws.Range("D16").Value = "'0.006026"
$ws.Range("E16").Value = "'-1.36%"
$ws.Range("D17").Value = "'3.459"
$ws.Range("D18").Value = "'3.166"
$ws.Range("E18").Value = "'-0.96%"
$ws.Range("D19").Value = "'2.166"
$ws.Range("E19").Value = "'-1.73%"
$ws.Range("D20").Value = "'0.3112"
$ws.Range("E20").Value = "'0.17%"
$ws.Range("D21").Value = "'0.1297"
$ws.Range("E21").Value = "'1.55%"
$ws.Range("D22").Value = "'4.139"
$ws.Range("E22").Value = "'7.47%"
$ws.Range("D23").Value = "'0.04237"
$ws.Range("E23").Value = "'0.40%"
$ws.Range("E24").Value = "'-0.60%"
$ws.Range("E25").Value = "'-6.20%"
$ws.Range("E26").Value = "'0.02%"
$ws.Range("D27").Value = "'0.0001715"
$ws.Range("E27").Value = "'13.54%"
$ws.Range("D40").Value = "'0.03841"
$ws.Range("E40").Value = "'0.12%"
$ws.Range("D41").Value = "'0.1113"
$ws.Range("E41").Value = "'0.92%"
$ws.Range("D42").Value = "'0.003841"
$ws.Range("E42").Value = "'-38.52%"
$ws.Range("D43").Value = "'0.002419"
$ws.Range("E43").Value = "'1.27%"
$ws.Range("D44").Value = "'0.01510"
$ws.Range("E44").Value = "'28.14%"
$ws.Range("D45").Value = "'0.00005114"
$ws.Range("E45").Value = "'-6.42%"
$ws.Range("E46").Value = "'-0.02%"
$ws.Range("E47").Value = "'-16.73%"
$ws.Range("E48").Value = "'-3.24%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'-0.02%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'-0.02%"
